$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 13554.667
$ws.Range("I32").Value = 998
$ws.Range("J32").Value = 15124.25
$ws.Range("K32").Value = 998
$ws.Range("L32").Value = 15124.25
$ws.Range("M32").Value = -672
$ws.Range("N32").Value = -15776.25
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
# Row 74
$ws.Range("H74").Value = 25697014
$ws.Range("I74").Value = 25697014
$ws.Range("K74").Value = 25697014
$ws.Range("M74").Value = -25696078
# Row 77
$ws.Range("H77").Value = 25697014
$ws.Range("I77").Value = 25697014
$ws.Range("K77").Value = 128485070
$ws.Range("M77").Value = -128480390
# Row 97
$ws.Range("H97").Value = 3511.75
$ws.Range("J97").Value = 3511.75
$ws.Range("L97").Value = 10535.25
$ws.Range("N97").Value = -11527.25
# Row 138
$ws.Range("H138").Value = 2530.0889
$ws.Range("I138").Value = 2334.6843
$ws.Range("K138").Value = 7004.0529
$ws.Range("M138").Value = -1864.0529

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5274.829
$ws.Range("I32").Value = 3558.625
$ws.Range("K32").Value = 3558.625
$ws.Range("M32").Value = -3271.625
# Row 74
$ws.Range("H74").Value = 71435880
$ws.Range("I74").Value = 100008530
$ws.Range("K74").Value = 100008530
$ws.Range("M74").Value = -100007656
# Row 77
$ws.Range("H77").Value = 71435880
$ws.Range("I77").Value = 100008530
$ws.Range("K77").Value = 500042650
$ws.Range("M77").Value = -500038282
# Row 132
$ws.Range("H132").Value = 8337893
$ws.Range("J132").Value = 3469
$ws.Range("L132").Value = 10407
$ws.Range("N132").Value = -15467
# Row 135
$ws.Range("H135").Value = 99932.336
$ws.Range("J135").Value = 99932.336
$ws.Range("L135").Value = 99932.336
$ws.Range("N135").Value = -110072.336

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 46817.863
$ws.Range("I107").Value = 1293.1177
$ws.Range("K107").Value = 1293.1177
$ws.Range("M107").Value = 626.8823
# Row 117
$ws.Range("H117").Value = 100000
$ws.Range("J117").Value = 100000
$ws.Range("L117").Value = 100000
$ws.Range("N117").Value = -109178

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 33599.332
$ws.Range("I22").Value = 33599.332
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 33599.332
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -33249.332
$ws.Range("N22").ClearContents()
# Row 31
$ws.Range("H31").Value = 8029.309
$ws.Range("I31").Value = 6086.276
$ws.Range("J31").Value = 10196.538
$ws.Range("K31").Value = 6086.276
$ws.Range("L31").Value = 10196.538
$ws.Range("M31").Value = -5791.276
$ws.Range("N31").Value = -10786.538
# Row 34
$ws.Range("H34").Value = 8029.309
$ws.Range("I34").Value = 6086.276
$ws.Range("J34").Value = 10196.538
$ws.Range("K34").Value = 6086.276
$ws.Range("L34").Value = 10196.538
$ws.Range("M34").Value = -5884.276
$ws.Range("N34").Value = -10600.538
# Row 58
$ws.Range("H58").Value = 25007244
$ws.Range("I58").Value = 33341604
$ws.Range("J58").Value = 4165.2
$ws.Range("K58").Value = 33341604
$ws.Range("L58").Value = 4165.2
$ws.Range("M58").Value = -33341401
$ws.Range("N58").Value = -4571.2
# Row 125
$ws.Range("H125").Value = 57500
$ws.Range("J125").Value = 57500
$ws.Range("L125").Value = 57500
$ws.Range("N125").Value = -62420
# Row 136
$ws.Range("H136").Value = 25007244
$ws.Range("I136").Value = 33341604
$ws.Range("J136").Value = 4165.2
$ws.Range("K136").Value = 100024812
$ws.Range("L136").Value = 12495.6
$ws.Range("M136").Value = -100022262
$ws.Range("N136").Value = -17595.6

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 97
$ws.Range("H97").Value = 1487.2858
$ws.Range("I97").Value = 92.22221999999999
$ws.Range("K97").Value = 276.66666
$ws.Range("M97").Value = 219.33334
# Row 119
$ws.Range("H119").Value = 6732
$ws.Range("I119").Value = 2309.6667
$ws.Range("J119").Value = 19999
$ws.Range("K119").Value = 6929.000100000001
$ws.Range("L119").Value = 59997
$ws.Range("M119").Value = -2091.000100000001
$ws.Range("N119").Value = -69673

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 789844.9
$ws.Range("J7").Value = 21666
$ws.Range("L7").Value = 21666
$ws.Range("N7").Value = -21890
# Row 8
$ws.Range("H8").Value = 789844.9
$ws.Range("J8").Value = 21666
$ws.Range("L8").Value = 21666
$ws.Range("N8").Value = -21944
# Row 102
$ws.Range("H102").Value = 4327.9287
$ws.Range("I102").Value = 814.6923
$ws.Range("K102").Value = 814.6923
$ws.Range("M102").Value = 807.3077
# Row 107
$ws.Range("H107").Value = 703
$ws.Range("I107").Value = 380.625
$ws.Range("K107").Value = 380.625
$ws.Range("M107").Value = 1539.375
# Row 124
$ws.Range("H124").Value = 73017
$ws.Range("J124").Value = 73017
$ws.Range("L124").Value = 73017
$ws.Range("N124").Value = -82837
# Row 132
$ws.Range("H132").Value = 7358409
$ws.Range("I132").Value = 8337296
$ws.Range("J132").Value = 16755
$ws.Range("K132").Value = 25011888
$ws.Range("L132").Value = 50265
$ws.Range("M132").Value = -25009358
$ws.Range("N132").Value = -55325

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3497.25
$ws.Range("I40").Value = 3497.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3497.25
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3361.25
$ws.Range("N40").ClearContents()
# Row 122
$ws.Range("H122").Value = 5498.3477
$ws.Range("I122").Value = 5472.9
$ws.Range("J122").Value = 5668
$ws.Range("K122").Value = 16418.7
$ws.Range("L122").Value = 17004
$ws.Range("M122").Value = -13968.7
$ws.Range("N122").Value = -21904
# Row 130
$ws.Range("H130").Value = 59999.5
$ws.Range("J130").Value = 59999.5
$ws.Range("L130").Value = 59999.5
$ws.Range("N130").Value = -70039.5
# Row 132
$ws.Range("H132").Value = 18471560
$ws.Range("I132").Value = 20010698
$ws.Range("K132").Value = 60032094
$ws.Range("M132").Value = -60029564
# Row 140
$ws.Range("H140").Value = 73553
$ws.Range("J140").Value = 73553
$ws.Range("L140").Value = 73553
$ws.Range("N140").Value = -83913

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 95000
$ws.Range("J46").Value = 95000
$ws.Range("L46").Value = 95000
$ws.Range("N46").Value = -95462
# Row 110
$ws.Range("H110").Value = 158249.5
$ws.Range("J110").Value = 158249.5
$ws.Range("L110").Value = 158249.5
$ws.Range("N110").Value = -166429.5
# Row 112
$ws.Range("H112").Value = 37866.2
$ws.Range("J112").Value = 37866.2
$ws.Range("L112").Value = 37866.2
$ws.Range("N112").Value = -40820.2
# Row 116
$ws.Range("H116").Value = 79000
$ws.Range("J116").Value = 79000
$ws.Range("L116").Value = 79000
$ws.Range("N116").Value = -88178
# Row 123
$ws.Range("H123").Value = 119999
$ws.Range("J123").Value = 119999
$ws.Range("L123").Value = 119999
$ws.Range("N123").Value = -129799
# Row 134
$ws.Range("H134").Value = 95000
$ws.Range("J134").Value = 95000
$ws.Range("L134").Value = 285000
$ws.Range("N134").Value = -290070
# Row 136
$ws.Range("H136").Value = 41667470
$ws.Range("I136").Value = 45454920
$ws.Range("K136").Value = 136364760
$ws.Range("M136").Value = -136362210
# Row 137
$ws.Range("H137").Value = 124000
$ws.Range("J137").Value = 124000
$ws.Range("L137").Value = 124000
$ws.Range("N137").Value = -134200
